# Append 1 row of data to the bottom of the used range (Sheet1!A11:H11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current used range.
$newRow = $ws.UsedRange.Rows.Count + 1

$rowValues = @(
    '',
    'يامن ',
    '23',
    'الصمود',
    'الرحلة 3',
    'C2',
    'NRC',
    '٠١‏/٠٥‏/٢٠٢٥ ٠٦:٥٥:٥٤ م'
)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $cell = $ws.Cells.Item($newRow, $i + 1)
    # Leading apostrophe forces text storage (no numeric/date coercion),
    # matching how the source data is stored as plain text.
    $cell.Value = "'" + $rowValues[$i]
    $cell.Style = "Normal"
}
